$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2304235
$ws.Range("J17").Value = 2304235
$ws.Range("L17").Value = 6912705
$ws.Range("N17").Value = -6913041
$ws.Range("H18").Value = 3262.6667
$ws.Range("I18").Value = 3262.6667
$ws.Range("K18").Value = 3262.6667
$ws.Range("M18").Value = -2978.6667
$ws.Range("H40").Value = 4649.0625
$ws.Range("J40").Value = 4341.143
$ws.Range("L40").Value = 4341.143
$ws.Range("N40").Value = -4691.143
$ws.Range("H58").Value = 568.7778
$ws.Range("J58").Value = 2016.5
$ws.Range("L58").Value = 6049.5
$ws.Range("N58").Value = -6349.5
$ws.Range("H64").Value = 10198.4
$ws.Range("I64").Value = 9347
$ws.Range("K64").Value = 9347
$ws.Range("M64").Value = -9099
$ws.Range("H67").Value = 10198.4
$ws.Range("I67").Value = 9347
$ws.Range("K67").Value = 9347
$ws.Range("M67").Value = -8489
$ws.Range("H69").Value = 53702.668
$ws.Range("J69").Value = 68054.5
$ws.Range("L69").Value = 204163.5
$ws.Range("N69").Value = -205911.5
$ws.Range("H72").Value = 53702.668
$ws.Range("J72").Value = 68054.5
$ws.Range("L72").Value = 612490.5
$ws.Range("N72").Value = -621226.5
$ws.Range("H80").Value = 69445060
$ws.Range("I80").Value = 76923460
$ws.Range("K80").Value = 230770380
$ws.Range("M80").Value = -230769382
$ws.Range("H83").Value = 69445060
$ws.Range("I83").Value = 76923460
$ws.Range("K83").Value = 692311140
$ws.Range("M83").Value = -692306148
$ws.Range("H86").Value = 78434170
$ws.Range("I86").Value = 76926060
$ws.Range("K86").Value = 76926060
$ws.Range("M86").Value = -76924937
$ws.Range("H88").Value = 1013.5
$ws.Range("I88").Value = 1220
$ws.Range("J88").Value = 961.875
$ws.Range("K88").Value = 1220
$ws.Range("L88").Value = 961.875
$ws.Range("M88").Value = -814
$ws.Range("N88").Value = -1773.875
$ws.Range("H89").Value = 78434170
$ws.Range("I89").Value = 76926060
$ws.Range("K89").Value = 384630300
$ws.Range("M89").Value = -384624684
$ws.Range("H91").Value = 1013.5
$ws.Range("I91").Value = 1220
$ws.Range("J91").Value = 961.875
$ws.Range("K91").Value = 1220
$ws.Range("L91").Value = 961.875
$ws.Range("M91").Value = 184
$ws.Range("N91").Value = -3769.875
$ws.Range("H103").Value = 2331.6667
$ws.Range("J103").Value = 1000
$ws.Range("L103").Value = 3000
$ws.Range("N103").Value = -4172
$ws.Range("H111").Value = 4351.8184
$ws.Range("I111").Value = 4208
$ws.Range("K111").Value = 12624
$ws.Range("M111").Value = -9557
$ws.Range("H113").Value = 7891.8335
$ws.Range("I113").Value = 6604
$ws.Range("J113").Value = 9694.799999999999
$ws.Range("K113").Value = 6604
$ws.Range("L113").Value = 9694.799999999999
$ws.Range("M113").Value = -3350
$ws.Range("N113").Value = -16202.8
$ws.Range("H125").Value = 77799800
$ws.Range("J125").Value = 4182.1665
$ws.Range("L125").Value = 37639.4985
$ws.Range("N125").Value = -42559.4985
$ws.Range("H132").Value = 38585.035
$ws.Range("I132").Value = 48385.59
$ws.Range("K132").Value = 145156.77
$ws.Range("M132").Value = -142626.77
$ws.Range("H137").Value = 4350061.5
$ws.Range("I137").Value = 1771.5294
$ws.Range("J137").Value = 16670217
$ws.Range("K137").Value = 5314.5882
$ws.Range("L137").Value = 50010651
$ws.Range("M137").Value = -2764.5882
$ws.Range("N137").Value = -50015751
$ws.Range("H138").Value = 8153.6665
$ws.Range("I138").Value = 12581.167
$ws.Range("J138").Value = 5939.9165
$ws.Range("K138").Value = 37743.501
$ws.Range("L138").Value = 17819.7495
$ws.Range("M138").Value = -32603.501
$ws.Range("N138").Value = -28099.7495
$ws.Range("H141").Value = 2098.5
$ws.Range("I141").Value = 2531.3333
$ws.Range("K141").Value = 7593.999899999999
$ws.Range("M141").Value = -2413.999899999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6370061.5
$ws.Range("I32").Value = 8380454
$ws.Range("K32").Value = 8380454
$ws.Range("M32").Value = -8380167
$ws.Range("H45").Value = 2756.6
$ws.Range("I45").Value = 2851.3333
$ws.Range("K45").Value = 2851.3333
$ws.Range("M45").Value = -2474.3333
$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()
$ws.Range("H61").Value = 2003153.9
$ws.Range("I61").Value = 12905.333
$ws.Range("K61").Value = 12905.333
$ws.Range("M61").Value = -12693.333
$ws.Range("H63").Value = 10763.375
$ws.Range("I63").Value = 4388.143
$ws.Range("K63").Value = 4388.143
$ws.Range("M63").Value = -3702.143
$ws.Range("H66").Value = 10763.375
$ws.Range("I66").Value = 4388.143
$ws.Range("K66").Value = 21940.715
$ws.Range("M66").Value = -18508.715
$ws.Range("H74").Value = 440505.12
$ws.Range("I74").Value = 3974.361
$ws.Range("K74").Value = 3974.361
$ws.Range("M74").Value = -3100.361
$ws.Range("H77").Value = 440505.12
$ws.Range("I77").Value = 3974.361
$ws.Range("K77").Value = 19871.805
$ws.Range("M77").Value = -15503.805
$ws.Range("H88").Value = 3037.182
$ws.Range("J88").Value = 3081.889
$ws.Range("L88").Value = 3081.889
$ws.Range("N88").Value = -3893.889
$ws.Range("H91").Value = 3037.182
$ws.Range("J91").Value = 3081.889
$ws.Range("L91").Value = 3081.889
$ws.Range("N91").Value = -5889.889
$ws.Range("H102").Value = 4164.2607
$ws.Range("I102").Value = 3814.3333
$ws.Range("J102").Value = 5424
$ws.Range("K102").Value = 3814.3333
$ws.Range("L102").Value = 5424
$ws.Range("M102").Value = -2192.3333
$ws.Range("N102").Value = -8668
$ws.Range("H122").Value = 1418.6923
$ws.Range("I122").Value = 1408.4546
$ws.Range("K122").Value = 4225.3638
$ws.Range("M122").Value = -1775.3638
$ws.Range("H136").Value = 2003153.9
$ws.Range("I136").Value = 12905.333
$ws.Range("K136").Value = 38715.999
$ws.Range("M136").Value = -36165.999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1238.7142
$ws.Range("I20").Value = 1154.2
$ws.Range("J20").Value = 1450
$ws.Range("K20").Value = 1154.2
$ws.Range("L20").Value = 1450
$ws.Range("M20").Value = -907.2
$ws.Range("N20").Value = -1944
$ws.Range("H58").Value = 0
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("N58").ClearContents()
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()
$ws.Range("H86").Value = 3812.25
$ws.Range("I86").Value = 2499.6667
$ws.Range("J86").Value = 4599.8
$ws.Range("K86").Value = 2499.6667
$ws.Range("L86").Value = 4599.8
$ws.Range("M86").Value = -1376.6667
$ws.Range("N86").Value = -6845.8
$ws.Range("H89").Value = 3812.25
$ws.Range("I89").Value = 2499.6667
$ws.Range("J89").Value = 4599.8
$ws.Range("K89").Value = 12498.3335
$ws.Range("L89").Value = 22999
$ws.Range("M89").Value = -6882.333500000001
$ws.Range("N89").Value = -34231
$ws.Range("H105").Value = 10411.632
$ws.Range("I105").Value = 12633.3
$ws.Range("J105").Value = 7943.1113
$ws.Range("K105").Value = 12633.3
$ws.Range("L105").Value = 7943.1113
$ws.Range("M105").Value = -10886.3
$ws.Range("N105").Value = -11437.1113
$ws.Range("H107").Value = 10491.857
$ws.Range("I107").Value = 12420.485
$ws.Range("J107").Value = 5670.2856
$ws.Range("K107").Value = 12420.485
$ws.Range("L107").Value = 5670.2856
$ws.Range("M107").Value = -10500.485
$ws.Range("N107").Value = -9510.285599999999
$ws.Range("H110").Value = 120000
$ws.Range("J110").Value = 120000
$ws.Range("L110").Value = 120000
$ws.Range("N110").Value = -128180
$ws.Range("H134").Value = 19174528
$ws.Range("I134").Value = 34636.79
$ws.Range("J134").Value = 64289988
$ws.Range("K134").Value = 103910.37
$ws.Range("L134").Value = 192869964
$ws.Range("M134").Value = -101375.37
$ws.Range("N134").Value = -192875034

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5074.225
$ws.Range("I31").Value = 3130.1765
$ws.Range("J31").Value = 6511.1304
$ws.Range("K31").Value = 3130.1765
$ws.Range("L31").Value = 6511.1304
$ws.Range("M31").Value = -2835.1765
$ws.Range("N31").Value = -7101.1304
$ws.Range("H34").Value = 5074.225
$ws.Range("I34").Value = 3130.1765
$ws.Range("J34").Value = 6511.1304
$ws.Range("K34").Value = 3130.1765
$ws.Range("L34").Value = 6511.1304
$ws.Range("M34").Value = -2928.1765
$ws.Range("N34").Value = -6915.1304
$ws.Range("H58").Value = 3330.5
$ws.Range("I58").Value = 2311
$ws.Range("J58").Value = 4641.2856
$ws.Range("K58").Value = 2311
$ws.Range("L58").Value = 4641.2856
$ws.Range("M58").Value = -2108
$ws.Range("N58").Value = -5047.2856
$ws.Range("H62").Value = 3947.8
$ws.Range("J62").Value = 4444
$ws.Range("L62").Value = 4444
$ws.Range("N62").Value = -5692
$ws.Range("H65").Value = 3947.8
$ws.Range("J65").Value = 4444
$ws.Range("L65").Value = 22220
$ws.Range("N65").Value = -28460
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("M68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("M71").ClearContents()
$ws.Range("H74").Value = 333396670
$ws.Range("J74").Value = 333396670
$ws.Range("L74").Value = 333396670
$ws.Range("N74").Value = -333398418
$ws.Range("H77").Value = 333396670
$ws.Range("J77").Value = 333396670
$ws.Range("L77").Value = 1000190010
$ws.Range("N77").Value = -1000198746
$ws.Range("H86").Value = 9472.115
$ws.Range("I86").Value = 5263.875
$ws.Range("J86").Value = 16205.3
$ws.Range("K86").Value = 5263.875
$ws.Range("L86").Value = 16205.3
$ws.Range("M86").Value = -4140.875
$ws.Range("N86").Value = -18451.3
$ws.Range("H89").Value = 9472.115
$ws.Range("I89").Value = 5263.875
$ws.Range("J89").Value = 16205.3
$ws.Range("K89").Value = 26319.375
$ws.Range("L89").Value = 81026.5
$ws.Range("M89").Value = -20703.375
$ws.Range("N89").Value = -92258.5
$ws.Range("H93").Value = 13427.556
$ws.Range("I93").Value = 11550
$ws.Range("K93").Value = 11550
$ws.Range("M93").Value = -9678
$ws.Range("H130").Value = 614995
$ws.Range("J130").Value = 614995
$ws.Range("L130").Value = 614995
$ws.Range("N130").Value = -625035
$ws.Range("H131").Value = 0
$ws.Range("I131").Value = 0
$ws.Range("K131").Value = 0
$ws.Range("M131").ClearContents()
$ws.Range("H132").Value = 18524184
$ws.Range("I132").Value = 8170.6665
$ws.Range("J132").Value = 22227386
$ws.Range("K132").Value = 24511.9995
$ws.Range("L132").Value = 66682158
$ws.Range("M132").Value = -21981.9995
$ws.Range("N132").Value = -66687218
$ws.Range("H134").Value = 2805.7715
$ws.Range("I134").Value = 2147.76
$ws.Range("K134").Value = 6443.280000000001
$ws.Range("M134").Value = -3908.280000000001
$ws.Range("H136").Value = 3330.5
$ws.Range("I136").Value = 2311
$ws.Range("J136").Value = 4641.2856
$ws.Range("K136").Value = 6933
$ws.Range("L136").Value = 13923.8568
$ws.Range("M136").Value = -4383
$ws.Range("N136").Value = -19023.8568
$ws.Range("H140").Value = 145494
$ws.Range("J140").Value = 145494
$ws.Range("L140").Value = 145494
$ws.Range("N140").Value = -155854
$ws.Range("H141").Value = 595277.2
$ws.Range("J141").Value = 586332.75
$ws.Range("L141").Value = 586332.75
$ws.Range("N141").Value = -596692.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 911.63635
$ws.Range("I2").Value = 430.33334
$ws.Range("J2").Value = 1160.5862
$ws.Range("K2").Value = 2582.00004
$ws.Range("L2").Value = 6963.5172
$ws.Range("M2").Value = -2469.00004
$ws.Range("N2").Value = -7189.5172
$ws.Range("H4").Value = 3018799.8
$ws.Range("I4").Value = 6864813
$ws.Range("K4").Value = 20594439
$ws.Range("M4").Value = -20594327
$ws.Range("H23").Value = 77062.08
$ws.Range("I23").Value = 256.6
$ws.Range("K23").Value = 769.8000000000001
$ws.Range("M23").Value = -534.8000000000001
$ws.Range("H50").Value = 6640.077
$ws.Range("I50").Value = 3930.8572
$ws.Range("K50").Value = 11792.5716
$ws.Range("M50").Value = -11311.5716
$ws.Range("H53").Value = 6640.077
$ws.Range("I53").Value = 3930.8572
$ws.Range("K53").Value = 11792.5716
$ws.Range("M53").Value = -11311.5716
$ws.Range("H113").Value = 1468.091
$ws.Range("I113").Value = 1447.3334
$ws.Range("K113").Value = 4342.0002
$ws.Range("M113").Value = -2172.0002
$ws.Range("H122").Value = 2447637.5
$ws.Range("J122").Value = 884637.7
$ws.Range("L122").Value = 7961739.3
$ws.Range("N122").Value = -7966639.3
$ws.Range("H132").Value = 1313
$ws.Range("I132").Value = 1196.6666
$ws.Range("K132").Value = 10769.9994
$ws.Range("M132").Value = -8239.999400000001
$ws.Range("H136").Value = 2304.8333
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H35").Value = 21999.5
$ws.Range("I35").Value = 21999.5
$ws.Range("K35").Value = 21999.5
$ws.Range("M35").Value = -21701.5
$ws.Range("H64").Value = 60000
$ws.Range("I64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("M64").ClearContents()
$ws.Range("H67").Value = 60000
$ws.Range("I67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("M67").ClearContents()
$ws.Range("H70").Value = 16809
$ws.Range("I70").Value = 15804.429
$ws.Range("K70").Value = 15804.429
$ws.Range("M70").Value = -15534.429
$ws.Range("H73").Value = 16809
$ws.Range("I73").Value = 15804.429
$ws.Range("K73").Value = 15804.429
$ws.Range("M73").Value = -14868.429
$ws.Range("H80").Value = 9449333
$ws.Range("J80").Value = 45665092
$ws.Range("L80").Value = 45665092
$ws.Range("N80").Value = -45667088
$ws.Range("H83").Value = 9449333
$ws.Range("J83").Value = 45665092
$ws.Range("L83").Value = 228325460
$ws.Range("N83").Value = -228335444
$ws.Range("H102").Value = 11905396
$ws.Range("I102").Value = 13889394
$ws.Range("K102").Value = 13889394
$ws.Range("M102").Value = -13887772
$ws.Range("H112").Value = 40000
$ws.Range("J112").Value = 40000
$ws.Range("L112").Value = 40000
$ws.Range("N112").Value = -42216
$ws.Range("H122").Value = 4196.8096
$ws.Range("I122").Value = 4996.5
$ws.Range("K122").Value = 14989.5
$ws.Range("M122").Value = -12539.5
$ws.Range("H132").Value = 3582974.2
$ws.Range("I132").Value = 2874.8076
$ws.Range("J132").Value = 11029581
$ws.Range("K132").Value = 8624.4228
$ws.Range("L132").Value = 33088743
$ws.Range("M132").Value = -6094.4228
$ws.Range("N132").Value = -33093803

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3589.9644
$ws.Range("I22").Value = 1681.5
$ws.Range("J22").Value = 5021.3125
$ws.Range("K22").Value = 1681.5
$ws.Range("L22").Value = 5021.3125
$ws.Range("M22").Value = -1386.5
$ws.Range("N22").Value = -5611.3125
$ws.Range("H27").Value = 3589.9644
$ws.Range("I27").Value = 1681.5
$ws.Range("J27").Value = 5021.3125
$ws.Range("K27").Value = 1681.5
$ws.Range("L27").Value = 5021.3125
$ws.Range("M27").Value = -1574.5
$ws.Range("N27").Value = -5235.3125
$ws.Range("H34").Value = 13888.75
$ws.Range("I34").Value = 10185
$ws.Range("J34").Value = 25000
$ws.Range("K34").Value = 10185
$ws.Range("L34").Value = 25000
$ws.Range("M34").Value = -10013
$ws.Range("N34").Value = -25344
$ws.Range("H42").Value = 270499.5
$ws.Range("I42").Value = 27333
$ws.Range("J42").Value = 999999
$ws.Range("K42").Value = 27333
$ws.Range("L42").Value = 999999
$ws.Range("M42").Value = -26770
$ws.Range("N42").Value = -1001125
$ws.Range("H46").Value = 12523.158
$ws.Range("J46").Value = 4998.3335
$ws.Range("L46").Value = 4998.3335
$ws.Range("N46").Value = -5374.3335
$ws.Range("H49").Value = 270499.5
$ws.Range("I49").Value = 27333
$ws.Range("J49").Value = 999999
$ws.Range("K49").Value = 27333
$ws.Range("L49").Value = 999999
$ws.Range("M49").Value = -27186
$ws.Range("N49").Value = -1000293
$ws.Range("H61").Value = 7145141
$ws.Range("I61").Value = 8697063
$ws.Range("K61").Value = 8697063
$ws.Range("M61").Value = -8696861
$ws.Range("H68").Value = 7113.375
$ws.Range("I68").Value = 5999.6665
$ws.Range("J68").Value = 7781.6
$ws.Range("K68").Value = 5999.6665
$ws.Range("L68").Value = 7781.6
$ws.Range("M68").Value = -5250.6665
$ws.Range("N68").Value = -9279.6
$ws.Range("H71").Value = 7113.375
$ws.Range("I71").Value = 5999.6665
$ws.Range("J71").Value = 7781.6
$ws.Range("K71").Value = 29998.3325
$ws.Range("L71").Value = 38908
$ws.Range("M71").Value = -26254.3325
$ws.Range("N71").Value = -46396
$ws.Range("H100").Value = 2453
$ws.Range("I100").Value = 1674.25
$ws.Range("J100").Value = 3699
$ws.Range("K100").Value = 1674.25
$ws.Range("L100").Value = 3699
$ws.Range("M100").Value = -1133.25
$ws.Range("N100").Value = -4781
$ws.Range("H113").Value = 7145141
$ws.Range("I113").Value = 8697063
$ws.Range("K113").Value = 8697063
$ws.Range("M113").Value = -8694893
$ws.Range("H122").Value = 4986.846
$ws.Range("I122").Value = 4466.8125
$ws.Range("K122").Value = 13400.4375
$ws.Range("M122").Value = -10950.4375
$ws.Range("H131").Value = 80000
$ws.Range("J131").Value = 80000
$ws.Range("L131").Value = 80000
$ws.Range("N131").Value = -90080
$ws.Range("H132").Value = 12086.238
$ws.Range("I132").Value = 5490.1113
$ws.Range("K132").Value = 16470.3339
$ws.Range("M132").Value = -13940.3339
$ws.Range("H133").Value = 170000
$ws.Range("J133").Value = 170000
$ws.Range("L133").Value = 170000
$ws.Range("N133").Value = -175060

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 0
$ws.Range("I54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("M54").ClearContents()
$ws.Range("H62").Value = 9498.5
$ws.Range("I62").Value = 9498.5
$ws.Range("K62").Value = 9498.5
$ws.Range("M62").Value = -8874.5
$ws.Range("H65").Value = 9498.5
$ws.Range("I65").Value = 9498.5
$ws.Range("K65").Value = 47492.5
$ws.Range("M65").Value = -44372.5
$ws.Range("H74").Value = 34052.5
$ws.Range("J74").Value = 22570
$ws.Range("L74").Value = 22570
$ws.Range("N74").Value = -24442
$ws.Range("H76").Value = 129999
$ws.Range("J76").Value = 129999
$ws.Range("L76").Value = 129999
$ws.Range("N76").Value = -130629
$ws.Range("H77").Value = 34052.5
$ws.Range("J77").Value = 22570
$ws.Range("L77").Value = 67710
$ws.Range("N77").Value = -77070
$ws.Range("H79").Value = 129999
$ws.Range("J79").Value = 129999
$ws.Range("L79").Value = 129999
$ws.Range("N79").Value = -132183
$ws.Range("H100").Value = 1768.909
$ws.Range("I100").Value = 1663.8334
$ws.Range("K100").Value = 3327.6668
$ws.Range("M100").Value = -2786.6668
$ws.Range("H107").Value = 66667356
$ws.Range("I107").Value = 713.55554
$ws.Range("J107").Value = 166667310
$ws.Range("K107").Value = 2140.66662
$ws.Range("L107").Value = 500001930
$ws.Range("M107").Value = -220.66662
$ws.Range("N107").Value = -500005770
$ws.Range("H132").Value = 47690.863
$ws.Range("I132").Value = 200335.8
$ws.Range("K132").Value = 601007.3999999999
$ws.Range("M132").Value = -598477.3999999999
$ws.Range("H136").Value = 9224.727999999999
$ws.Range("I136").Value = 2353.75
$ws.Range("K136").Value = 7061.25
$ws.Range("M136").Value = -4511.25
$ws.Range("H138").Value = 0
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("M138").ClearContents()
$ws.Range("N138").ClearContents()
